$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 36,11
$arr[0,0] = "Hydrogen"
$arr[0,1] = 2030
$arr[0,2] = ""
$arr[0,3] = ""
$arr[0,4] = ""
$arr[0,5] = [double]"0.000203908691209593"
$arr[0,6] = ""
$arr[0,7] = [double]"3.133351628447623e-10"
$arr[0,8] = [double]"4.024665539268523e-05"
$arr[0,9] = ""
$arr[0,10] = ""
$arr[1,0] = "Methanol"
$arr[1,1] = 2030
$arr[1,2] = ""
$arr[1,3] = ""
$arr[1,4] = ""
$arr[1,5] = ""
$arr[1,6] = ""
$arr[1,7] = ""
$arr[1,8] = ""
$arr[1,9] = ""
$arr[1,10] = ""
$arr[2,0] = "Ammonia"
$arr[2,1] = 2030
$arr[2,2] = ""
$arr[2,3] = [double]"0.01403248634993889"
$arr[2,4] = ""
$arr[2,5] = ""
$arr[2,6] = ""
$arr[2,7] = ""
$arr[2,8] = ""
$arr[2,9] = ""
$arr[2,10] = ""
$arr[3,0] = "Synthetic Gases"
$arr[3,1] = 2030
$arr[3,2] = ""
$arr[3,3] = ""
$arr[3,4] = ""
$arr[3,5] = ""
$arr[3,6] = ""
$arr[3,7] = ""
$arr[3,8] = ""
$arr[3,9] = ""
$arr[3,10] = ""
$arr[4,0] = "Biogenic Gases"
$arr[4,1] = 2030
$arr[4,2] = ""
$arr[4,3] = ""
$arr[4,4] = [double]"9.827902782993092e-05"
$arr[4,5] = [double]"4.761502652839609e-05"
$arr[4,6] = ""
$arr[4,7] = ""
$arr[4,8] = [double]"1.010880051879498e-05"
$arr[4,9] = ""
$arr[4,10] = ""
$arr[5,0] = "Fossil Gases"
$arr[5,1] = 2030
$arr[5,2] = ""
$arr[5,3] = ""
$arr[5,4] = ""
$arr[5,5] = [double]"0.0007846670604052373"
$arr[5,6] = ""
$arr[5,7] = ""
$arr[5,8] = [double]"3.445060087347233e-05"
$arr[5,9] = ""
$arr[5,10] = ""
$arr[6,0] = "Synthetic Liquids"
$arr[6,1] = 2030
$arr[6,2] = ""
$arr[6,3] = ""
$arr[6,4] = ""
$arr[6,5] = ""
$arr[6,6] = ""
$arr[6,7] = ""
$arr[6,8] = ""
$arr[6,9] = ""
$arr[6,10] = ""
$arr[7,0] = "Biogenic Liquids"
$arr[7,1] = 2030
$arr[7,2] = ""
$arr[7,3] = ""
$arr[7,4] = ""
$arr[7,5] = [double]"0.003487021553041186"
$arr[7,6] = [double]"4.448074282126064e-06"
$arr[7,7] = [double]"0.0008573152441552"
$arr[7,8] = [double]"0.0015915965913147"
$arr[7,9] = [double]"0.0003263595920562"
$arr[7,10] = [double]"0.003283197281710558"
$arr[8,0] = "Fossil Liquids"
$arr[8,1] = 2030
$arr[8,2] = ""
$arr[8,3] = ""
$arr[8,4] = ""
$arr[8,5] = [double]"0.0360089500660149"
$arr[8,6] = [double]"3.132648279294923e-05"
$arr[8,7] = [double]"0.0078353990726053"
$arr[8,8] = [double]"0.0102817325873576"
$arr[8,9] = [double]"0.0019754552973177"
$arr[8,10] = [double]"0.03213083429072289"
$arr[9,0] = "Biomass [Solid]"
$arr[9,1] = 2030
$arr[9,2] = ""
$arr[9,3] = ""
$arr[9,4] = [double]"0.000471447495782198"
$arr[9,5] = ""
$arr[9,6] = ""
$arr[9,7] = ""
$arr[9,8] = ""
$arr[9,9] = ""
$arr[9,10] = ""
$arr[10,0] = "Renewable Energy Carrier"
$arr[10,1] = 2030
$arr[10,2] = ""
$arr[10,3] = ""
$arr[10,4] = [double]"0.000393890641552741"
$arr[10,5] = ""
$arr[10,6] = ""
$arr[10,7] = ""
$arr[10,8] = ""
$arr[10,9] = ""
$arr[10,10] = ""
$arr[11,0] = "Overall Demand"
$arr[11,1] = 2030
$arr[11,2] = ""
$arr[11,3] = [double]"0.01403248634993889"
$arr[11,4] = [double]"0.0009636171651648699"
$arr[11,5] = [double]"0.04053216239719931"
$arr[11,6] = [double]"3.577455707507529e-05"
$arr[11,7] = [double]"0.008692714630095664"
$arr[11,8] = [double]"0.01195813523545725"
$arr[11,9] = [double]"0.0023018148893739"
$arr[11,10] = [double]"0.03541403157243345"
$arr[12,0] = "Hydrogen"
$arr[12,1] = 2040
$arr[12,2] = ""
$arr[12,3] = ""
$arr[12,4] = ""
$arr[12,5] = [double]"0.0009889928412376682"
$arr[12,6] = ""
$arr[12,7] = [double]"2.622961394133191e-08"
$arr[12,8] = [double]"6.312641897502194e-05"
$arr[12,9] = ""
$arr[12,10] = ""
$arr[13,0] = "Methanol"
$arr[13,1] = 2040
$arr[13,2] = ""
$arr[13,3] = ""
$arr[13,4] = ""
$arr[13,5] = ""
$arr[13,6] = ""
$arr[13,7] = ""
$arr[13,8] = ""
$arr[13,9] = ""
$arr[13,10] = ""
$arr[14,0] = "Ammonia"
$arr[14,1] = 2040
$arr[14,2] = ""
$arr[14,3] = [double]"0.01292083825998528"
$arr[14,4] = ""
$arr[14,5] = ""
$arr[14,6] = ""
$arr[14,7] = ""
$arr[14,8] = ""
$arr[14,9] = ""
$arr[14,10] = ""
$arr[15,0] = "Synthetic Gases"
$arr[15,1] = 2040
$arr[15,2] = ""
$arr[15,3] = ""
$arr[15,4] = ""
$arr[15,5] = [double]"3.467063370201273e-10"
$arr[15,6] = ""
$arr[15,7] = ""
$arr[15,8] = [double]"2.152915032898363e-11"
$arr[15,9] = ""
$arr[15,10] = ""
$arr[16,0] = "Biogenic Gases"
$arr[16,1] = 2040
$arr[16,2] = ""
$arr[16,3] = ""
$arr[16,4] = [double]"0.0004169526670065891"
$arr[16,5] = [double]"6.248113136027162e-05"
$arr[16,6] = ""
$arr[16,7] = ""
$arr[16,8] = [double]"1.486320377195297e-05"
$arr[16,9] = ""
$arr[16,10] = ""
$arr[17,0] = "Fossil Gases"
$arr[17,1] = 2040
$arr[17,2] = ""
$arr[17,3] = ""
$arr[17,4] = ""
$arr[17,5] = [double]"0.0004126675344026147"
$arr[17,6] = ""
$arr[17,7] = ""
$arr[17,8] = [double]"3.717762580518311e-05"
$arr[17,9] = ""
$arr[17,10] = ""
$arr[18,0] = "Synthetic Liquids"
$arr[18,1] = 2040
$arr[18,2] = ""
$arr[18,3] = ""
$arr[18,4] = ""
$arr[18,5] = ""
$arr[18,6] = ""
$arr[18,7] = ""
$arr[18,8] = ""
$arr[18,9] = ""
$arr[18,10] = ""
$arr[19,0] = "Biogenic Liquids"
$arr[19,1] = 2040
$arr[19,2] = ""
$arr[19,3] = ""
$arr[19,4] = ""
$arr[19,5] = [double]"0.001389708746792253"
$arr[19,6] = [double]"7.251094900161178e-06"
$arr[19,7] = [double]"0.0010436587917408"
$arr[19,8] = [double]"0.0010829734611768"
$arr[19,9] = [double]"0.0003964739367207"
$arr[19,10] = [double]"0.003716633516751387"
$arr[20,0] = "Fossil Liquids"
$arr[20,1] = 2040
$arr[20,2] = ""
$arr[20,3] = ""
$arr[20,4] = ""
$arr[20,5] = [double]"0.0091903822788067"
$arr[20,6] = [double]"3.367771066565236e-05"
$arr[20,7] = [double]"0.0073873172301697"
$arr[20,8] = [double]"0.004741817687910699"
$arr[20,9] = [double]"0.001753039842156"
$arr[20,10] = [double]"0.03116482720734355"
$arr[21,0] = "Biomass [Solid]"
$arr[21,1] = 2040
$arr[21,2] = ""
$arr[21,3] = ""
$arr[21,4] = [double]"0.0004564710695266859"
$arr[21,5] = ""
$arr[21,6] = ""
$arr[21,7] = ""
$arr[21,8] = ""
$arr[21,9] = ""
$arr[21,10] = ""
$arr[22,0] = "Renewable Energy Carrier"
$arr[22,1] = 2040
$arr[22,2] = ""
$arr[22,3] = ""
$arr[22,4] = [double]"0.001560615532668226"
$arr[22,5] = ""
$arr[22,6] = ""
$arr[22,7] = ""
$arr[22,8] = ""
$arr[22,9] = ""
$arr[22,10] = ""
$arr[23,0] = "Overall Demand"
$arr[23,1] = 2040
$arr[23,2] = ""
$arr[23,3] = [double]"0.01292083825998528"
$arr[23,4] = [double]"0.002434039269201501"
$arr[23,5] = [double]"0.01204423287930585"
$arr[23,6] = [double]"4.092880556581354e-05"
$arr[23,7] = [double]"0.008431002251524442"
$arr[23,8] = [double]"0.005939958419168808"
$arr[23,9] = [double]"0.0021495137788767"
$arr[23,10] = [double]"0.03488146072409494"
$arr[24,0] = "Hydrogen"
$arr[24,1] = 2050
$arr[24,2] = ""
$arr[24,3] = ""
$arr[24,4] = ""
$arr[24,5] = [double]"0.001367417416914807"
$arr[24,6] = ""
$arr[24,7] = [double]"4.445708287726924e-08"
$arr[24,8] = [double]"0.0001016994180446112"
$arr[24,9] = ""
$arr[24,10] = ""
$arr[25,0] = "Methanol"
$arr[25,1] = 2050
$arr[25,2] = ""
$arr[25,3] = ""
$arr[25,4] = ""
$arr[25,5] = ""
$arr[25,6] = ""
$arr[25,7] = ""
$arr[25,8] = ""
$arr[25,9] = ""
$arr[25,10] = ""
$arr[26,0] = "Ammonia"
$arr[26,1] = 2050
$arr[26,2] = ""
$arr[26,3] = [double]"0.01300363302909579"
$arr[26,4] = ""
$arr[26,5] = ""
$arr[26,6] = ""
$arr[26,7] = ""
$arr[26,8] = ""
$arr[26,9] = ""
$arr[26,10] = ""
$arr[27,0] = "Synthetic Gases"
$arr[27,1] = 2050
$arr[27,2] = ""
$arr[27,3] = ""
$arr[27,4] = ""
$arr[27,5] = [double]"1.918199021657344e-09"
$arr[27,6] = ""
$arr[27,7] = ""
$arr[27,8] = [double]"7.04811139942833e-10"
$arr[27,9] = ""
$arr[27,10] = ""
$arr[28,0] = "Biogenic Gases"
$arr[28,1] = 2050
$arr[28,2] = ""
$arr[28,3] = ""
$arr[28,4] = [double]"0.001038000926619892"
$arr[28,5] = [double]"8.924190680993359e-06"
$arr[28,6] = ""
$arr[28,7] = ""
$arr[28,8] = [double]"4.08516418395332e-06"
$arr[28,9] = ""
$arr[28,10] = ""
$arr[29,0] = "Fossil Gases"
$arr[29,1] = 2050
$arr[29,2] = ""
$arr[29,3] = ""
$arr[29,4] = ""
$arr[29,5] = [double]"1.862511076346868e-05"
$arr[29,6] = ""
$arr[29,7] = ""
$arr[29,8] = [double]"1.291985523102099e-05"
$arr[29,9] = ""
$arr[29,10] = ""
$arr[30,0] = "Synthetic Liquids"
$arr[30,1] = 2050
$arr[30,2] = ""
$arr[30,3] = ""
$arr[30,4] = ""
$arr[30,5] = [double]"4.213050100346972e-12"
$arr[30,6] = [double]"2.670186332968756e-13"
$arr[30,7] = [double]"4.020933736636608e-11"
$arr[30,8] = [double]"1.686458052879049e-11"
$arr[30,9] = [double]"2.715380889178418e-12"
$arr[30,10] = [double]"2.55712462674254e-10"
$arr[31,0] = "Biogenic Liquids"
$arr[31,1] = 2050
$arr[31,2] = ""
$arr[31,3] = ""
$arr[31,4] = ""
$arr[31,5] = [double]"8.375249519702783e-05"
$arr[31,6] = [double]"1.295643350151406e-05"
$arr[31,7] = [double]"0.0013669681580208"
$arr[31,8] = [double]"0.0002817795873646104"
$arr[31,9] = [double]"0.0005099282662934"
$arr[31,10] = [double]"0.005293938082106364"
$arr[32,0] = "Fossil Liquids"
$arr[32,1] = 2050
$arr[32,2] = ""
$arr[32,3] = ""
$arr[32,4] = ""
$arr[32,5] = [double]"0.0002889457020113435"
$arr[32,6] = [double]"3.037586875746204e-05"
$arr[32,7] = [double]"0.006699345703858"
$arr[32,8] = [double]"0.0008284153699908"
$arr[32,9] = [double]"0.001508930617589"
$arr[32,10] = [double]"0.02904665318105566"
$arr[33,0] = "Biomass [Solid]"
$arr[33,1] = 2050
$arr[33,2] = ""
$arr[33,3] = ""
$arr[33,4] = [double]"0.0004575611005498155"
$arr[33,5] = ""
$arr[33,6] = ""
$arr[33,7] = ""
$arr[33,8] = ""
$arr[33,9] = ""
$arr[33,10] = ""
$arr[34,0] = "Renewable Energy Carrier"
$arr[34,1] = 2050
$arr[34,2] = ""
$arr[34,3] = ""
$arr[34,4] = [double]"0.003904132150802423"
$arr[34,5] = ""
$arr[34,6] = ""
$arr[34,7] = ""
$arr[34,8] = ""
$arr[34,9] = ""
$arr[34,10] = ""
$arr[35,0] = "Overall Demand"
$arr[35,1] = 2050
$arr[35,2] = ""
$arr[35,3] = [double]"0.01300363302909579"
$arr[35,4] = [double]"0.00539969417797213"
$arr[35,5] = [double]"0.001767666837979713"
$arr[35,6] = [double]"4.333230252599473e-05"
$arr[35,7] = [double]"0.008066358359171015"
$arr[35,8] = [double]"0.001228900116490717"
$arr[35,9] = [double]"0.002018858886597781"
$arr[35,10] = [double]"0.03434059151887449"

$ws.Range("A2:K37").Value = $arr

Write-Host "Done"
